# Updated to recent database structure / fixed interface
# - Replace the generic "from CEA, costs in USD-2015" reference note with a
#   more specific one wherever it is used.
# - Fill in the previously-blank cost/assumption columns (D, G) for the rows
#   that referenced that note, on both the HEATING and ELECTRICITY sheets.
# - Update sheet selections and the active sheet/tab (ELECTRICITY becomes the
#   active sheet instead of HEATING).

$wb = $excel.ActiveWorkbook

$newReference = "from CEA, costs in USD-2015, except for PEN and CO2, rest are assumptions"

# ---------------------------------------------------------------------
# HEATING sheet
# ---------------------------------------------------------------------
$wsHeating = $wb.Worksheets.Item("HEATING")

$wsHeating.Range("D5").Value = 0.8
$wsHeating.Range("G5").Value = 0.22
$wsHeating.Range("H5").Value = $newReference

$wsHeating.Range("D6").Value = 0.8
$wsHeating.Range("G6").Value = 0.22
$wsHeating.Range("H6").Value = $newReference

$wsHeating.Range("D7").Value = 0.8
$wsHeating.Range("G7").Value = 0.22
$wsHeating.Range("H7").Value = $newReference

# ---------------------------------------------------------------------
# ELECTRICITY sheet
# ---------------------------------------------------------------------
$wsElectricity = $wb.Worksheets.Item("ELECTRICITY")

$wsElectricity.Range("D5").Value = 0.99
$wsElectricity.Range("G5").Formula = "=0.22*0.75"
$wsElectricity.Range("H5").Value = $newReference

$wsElectricity.Range("D6").Value = 0.99
$wsElectricity.Range("G6").Formula = "=0.22*0.75"
$wsElectricity.Range("H6").Value = $newReference

$wsElectricity.Range("D7").Value = 0.99
$wsElectricity.Range("G7").Formula = "=0.22*0.75"
$wsElectricity.Range("H7").Value = $newReference

$wsElectricity.Range("D8").Value = 0.99
$wsElectricity.Range("G8").Formula = "=0.22*0.75"
$wsElectricity.Range("H8").Value = $newReference

# ---------------------------------------------------------------------
# Selections: HEATING keeps a selection but is no longer the active sheet;
# ELECTRICITY becomes the active sheet with its own selection.
# ---------------------------------------------------------------------
[void]$wsHeating.Range("G12").Select()

[void]$wsElectricity.Activate()
[void]$wsElectricity.Range("E14").Select()
